$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose displayed/stored text values change per the source diff.
# Price/volume columns hold plain text (not numeric) values, so each cell
# number format is set to Text ("@") before assignment - this stops Excel
# from reinterpreting strings like "314.85" or "2.24%" as numbers, which
# would silently alter precision / drop trailing zeros / drop the "%" sign.
$updates = [ordered]@{
    "D2" = "314.85"
    "E2" = "2.24%"
    "E3" = "2.30%"
    "D4" = "5.129"
    "E4" = "0.60%"
    "D5" = "0.08188"
    "E5" = "0.89%"
    "D6" = "1.980"
    "E6" = "0.58%"
    "D7" = "8.171"
    "E7" = "2.77%"
    "E8" = "-0.17%"
    "D9" = "0.1411"
    "E9" = "-1.06%"
    "D10" = "0.1958"
    "E10" = "-0.21%"
    "D11" = "0.09015"
    "E11" = "-0.90%"
    "D12" = "0.03494"
    "E12" = "-0.46%"
    "D13" = "0.09813"
    "E13" = "-0.17%"
    "D14" = "0.001401"
    "E14" = "-0.89%"
    "D15" = "0.005884"
    "E15" = "-3.66%"
    "D16" = "3.678"
    "E16" = "-0.07%"
    "D17" = "4.227"
    "E17" = "0.46%"
    "D18" = "3.174"
    "E18" = "-8.65%"
    "D19" = "0.3461"
    "E19" = "-0.04%"
    "D20" = "0.1351"
    "E20" = "3.66%"
    "D21" = "4.642"
    "E21" = "-3.64%"
    "D22" = "0.2421"
    "E22" = "-1.41%"
    "D23" = "0.04373"
    "E23" = "-1.09%"
    "D24" = "0.001244"
    "E24" = "1.97%"
    "D25" = "0.004803"
    "E25" = "-0.65%"
    "D27" = "0.0003994"
    "E27" = "-10.20%"
    "D39" = "0.02155"
    "E39" = "2.80%"
    "D40" = "0.05210"
    "E40" = "1.52%"
    "D41" = "0.007557"
    "E41" = "1.27%"
    "D42" = "0.009751"
    "E42" = "-3.83%"
    "D43" = "0.1374"
    "E43" = "1.07%"
    "D44" = "0.002114"
    "D45" = "0.009856"
    "E45" = "6.64%"
    "D46" = "0.00006370"
    "E46" = "1.85%"
    "E47" = "-0.30%"
    "E48" = "-9.97%"
    "D50" = "0.00002097"
    "E50" = "-0.30%"
    "E51" = "-0.30%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

